# Generate Report for Handoff
# Updates the localization-status report after a new handoff for b.md:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   - zh-cn / de-de sheets get fresh handoff file names + handoff datetimes for b.md
#   - de-de's b.md content-duplicate flag clears and both locales surface an
#     "out of date handback" error detail message
#   - Column P (Error Detail) is widened to fit the new message

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"
$overviewDate = "2016-08-27 18:45:18"

$zhHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhHandoffDate = "2016-08-27 18:45:14"

$deHandoffFile = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$deHandoffDate = "2016-08-27 18:45:18"

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cbdbf6949f58178faa2efdb33df005f3baa76ea3/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/24f9c75d1a90ce34c9554f8db94ca84671b06e02/e2e/b.md."

# ---- Overview sheet: row for b.md (row 3) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = $statusReady
$wsOverview.Range("F3").Value = $statusReady
$wsOverview.Range("G3").Value = $overviewDate

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
# row 2 = a.md, row 3 = b.md
$wsZh.Range("C2").Value = $statusReady
$wsZh.Range("C3").Value = $statusReady
# Leading apostrophe keeps this a text "False" (matching the other Content
# Duplicate cells) instead of auto-converting to a Boolean.
$wsZh.Range("F3").Value = "'False"
$wsZh.Range("G3").Value = $zhHandoffFile
$wsZh.Range("H3").Value = $zhHandoffDate
$wsZh.Range("P3").Value = $errorDetail
# ColumnWidth is in "characters"; the saved OOXML <col width> includes the
# ~0.8333 padding offset, so back that off to land on a stored width of 40.
$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusReady
$wsDe.Range("C3").Value = $statusReady
$wsDe.Range("F3").Value = "'False"
$wsDe.Range("G3").Value = $deHandoffFile
$wsDe.Range("H3").Value = $deHandoffDate
$wsDe.Range("P3").Value = $errorDetail
$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
